$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: instagram / https://www.instagram.com/pine
$ws.Range("A4").Value = "instagram"
$ws.Range("B4").Value = "https://www.instagram.com/pine"

# Turn B4 into a hyperlink, matching the style used by B2/B3
$ws.Hyperlinks.Add($ws.Range("B4"), "https://www.instagram.com/pine")
$ws.Range("B4").Style = "Hyperlink"

# Update the active selection to D11 (per diff)
$ws.Range("D11").Select()

Write-Host "done"
